# Generate Report for Handoff
# Updates the "Overview", "zh-cn" and "de-de" sheets to reflect that the
# 64799798-... file is now "Ready for handoff" (instead of "Handed back: in
# sync with en-US"), with refreshed handoff timestamps, and removes the
# second data row (the 6f99df9d-... file entry) from every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-23 04:54:10"

# Drop the 6f99df9d-... row (row 3) entirely, including its hyperlink.
$wsOverview.Range("A3").Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest
# Handoff File | Latest Handoff Datetime | Latest Target File | Latest
# Handback File | Latest Handback DateTime | ... | Handoff Reason | ...
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-23 04:54:07"

# Drop the 6f99df9d-... row (row 3) entirely, including its hyperlinks.
$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Range("D3").Hyperlinks.Delete()
$wsZhCn.Range("F3").Hyperlinks.Delete()
$wsZhCn.Range("G3").Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-23 04:54:10"

# Drop the 6f99df9d-... row (row 3) entirely, including its hyperlinks.
$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Range("D3").Hyperlinks.Delete()
$wsDeDe.Range("F3").Hyperlinks.Delete()
$wsDeDe.Range("G3").Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()
